$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.510.30"
$ws.Range("E2").Value = "  -2.20%  "

# Row 3
$ws.Range("D3").Value = "2.492.32"
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.80%  "

# Row 7
$ws.Range("E7").Value = "  -2.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.32%  "

# Row 11
$ws.Range("E11").Value = "  -1.69%  "

# Row 12
$ws.Range("E12").Value = "  +0.39%  "

# Row 13
$ws.Range("D13").Value = "2.876.76"
$ws.Range("E13").Value = "  -0.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.55%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.499.41"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("D18").Value = "41.473.33"
$ws.Range("E18").Value = "  -2.29%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.77%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0931"
$ws.Range("E20").Value = "  -0.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.38%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("E26").Value = "  -4.47%  "

# Row 27
$ws.Range("E27").Value = "  -3.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "

# Row 29
$ws.Range("E29").Value = "  -0.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.43%  "

# Row 33
$ws.Range("E33").Value = "  -2.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0756"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.81%  "

# Row 36
$ws.Range("E36").Value = "  -2.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.94%  "

# Row 38
$ws.Range("E38").Value = "  -3.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.92%  "

# Row 40
$ws.Range("E40").Value = "  -6.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

# Row 42
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.87%  "

# Row 44
$ws.Range("D44").Value = "1.992.34"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("E45").Value = "  -3.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.17%  "

# Row 48
$ws.Range("D48").Value = "2.739.83"
$ws.Range("E48").Value = "  -0.59%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.73%  "

# Row 51
$ws.Range("E51").Value = "  -4.53%  "
